# Generate Report for Archive
#
# The localization status report had its "Status" value updated from
# "Ready for handoff" to "In Translation" for every localized-file row,
# on every sheet that tracks it (Overview, zh-cn, de-de). After the
# shorter status text was written, the affected status columns were
# narrowed to fit the new (shorter) content.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update every cell currently showing the old status text -------------
# NOTE: cell values come back from Value2 with their native type (e.g. a
# genuine [bool] for TRUE/FALSE cells), and PowerShell's `-eq` coerces its
# right-hand operand to the type of the left-hand operand. Keeping the
# string literal on the left keeps the comparison a plain string compare
# regardless of the cell's actual type.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value2 = $newStatus
        }
    }
}

# --- Narrow the status columns to fit the shorter text --------------------
# Overview sheet: zh-cn (E) and de-de (F) status columns
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

# zh-cn sheet: Status column (C)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

# de-de sheet: Status column (C)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
